# Doing Updates for Financials
# Insert a new column before column D (pushing the existing D:K data to E:L)
# and populate the new column D with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; existing D:K columns shift right to E:L.
$ws.Columns("D:D").Insert()

# Copy number/date formatting from the (now shifted) column E into the new
# column D so the new cells keep the same look (date format for row 7/38/80,
# plain numeric for the rest) as the columns around them.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new period's values.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 816200
$ws.Range("D9").Value2 = 717600
$ws.Range("D10").Value2 = 98600
$ws.Range("D12").Value2 = 7300
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 1700
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 800100
$ws.Range("D18").Value2 = 16100
$ws.Range("D20").Value2 = 2200
$ws.Range("D21").Value2 = 28700
$ws.Range("D22").Value2 = 1100
$ws.Range("D23").Value2 = 17300
$ws.Range("D24").Value2 = 2300
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 15000
$ws.Range("D27").Value2 = 15000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -2200
$ws.Range("D33").Value2 = 15000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 15000
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 27400
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 142800
$ws.Range("D44").Value2 = 70000
$ws.Range("D45").Value2 = 5100
$ws.Range("D46").Value2 = 245300
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 56600
$ws.Range("D49").Value2 = 42400
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 9500
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 353800
$ws.Range("D57").Value2 = 76400
$ws.Range("D58").Value2 = 100
$ws.Range("D59").Value2 = 61600
$ws.Range("D60").Value2 = 138100
$ws.Range("D61").Value2 = 25500
$ws.Range("D62").Value2 = 4100
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 167000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 103600
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 186700
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 15000
$ws.Range("D83").Value2 = 10400
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 8000
$ws.Range("D91").Value2 = -9000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -14200
$ws.Range("D96").Value2 = -3500
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 100
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -6100

$wb.Save()
